$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.249.92"
$ws.Range("E2").Value = "  +0.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.906.56"
$ws.Range("E3").Value = "  +0.05%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "307.53"
$ws.Range("E5").Value = "  +0.36%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("E8").Value = "  +1.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07292"
$ws.Range("E9").Value = "  +0.56%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.66"
$ws.Range("E10").Value = "  +2.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.9041"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08175"
$ws.Range("E12").Value = "  -3.97%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.28"
$ws.Range("E13").Value = "  -0.48%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.367"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.531.05"
$ws.Range("E15").Value = "  -19.65%  "
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008677"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "14.76"
$ws.Range("E18").Value = "  +1.40%  "
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.281.12"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.119"
$ws.Range("E21").Value = "  +0.75%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.82"
$ws.Range("E22").Value = "  +1.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.508"
$ws.Range("E23").Value = "  +1.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.347"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.87"
$ws.Range("E25").Value = "  +1.88%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.737"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "116.72"
$ws.Range("E28").Value = "  +1.38%  "
$ws.Range("E29").Value = "  +0.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.859"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09247"
$ws.Range("E31").Value = "  -0.64%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8357"
$ws.Range("E32").Value = "  +3.75%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05070"
$ws.Range("E33").Value = "  +0.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.231"
$ws.Range("E34").Value = "  -1.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.986"
$ws.Range("E35").Value = "  +1.21%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.357"
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.723"
$ws.Range("E37").Value = "  +4.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5766"
$ws.Range("E38").Value = "  +0.81%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02005"
$ws.Range("E39").Value = "  +0.21%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.080"
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "9.110"
$ws.Range("E41").Value = "  -0.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.620"
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "116.22"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4910"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.22"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.645"
$ws.Range("E48").Value = "  +1.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "38.82"
$ws.Range("E49").Value = "  +3.19%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "64.43"
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").Value = "  +1.71%  "
